$wb = $excel.ActiveWorkbook

# Sheets 1-4: "Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
#             "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)"
# Same structural changes to each: add A1 header, strip style from A2:A12,
# and fix a few accented labels.
for ($i = 1; $i -le 4; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Add column header in A1, matching the style already used by B1 (bold,
    # bordered, centered) by copying its format.
    $ws.Range("A1").Value = "Fonte/Tecnologia"
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

    # Remove the bold/border style from the row-label cells A2:A12 (keep text).
    $ws.Range("A2:A12").ClearFormats()

    # Fix accented text in row labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."
}

$excel.CutCopyMode = $false

# Sheet 5: "Emissoes Totais (MtCO2eq)"
$ws5 = $wb.Worksheets.Item(5)

$ws5.Range("A1").Value = "Período"
$ws5.Range("B1").Copy() | Out-Null
$ws5.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

$ws5.Range("A2:A3").ClearFormats()
$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"

# Remove row 4 ("Teto") entirely.
$ws5.Range("A4:E4").Delete() | Out-Null

# Sheet 6: "Custo Total (bilhões de R$)"
$ws6 = $wb.Worksheets.Item(6)

$ws6.Range("A1").Value = "Tipo Expansão"
$ws6.Range("B1").Copy() | Out-Null
$ws6.Range("A1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Set B1's text to "2015" (matching the other sheets' header row) while
# keeping it a proper text cell (not a number) by copying a like cell.
$wb.Worksheets.Item(1).Range("B1").Copy() | Out-Null
$ws6.Range("B1").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = $false

$ws6.Range("A2:A3").ClearFormats()
$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("B2").Value = 556
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("B3").Value = 99
